# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" header on the Weekly Quantity sheet to
#    "Weekly_PO_Qty" and on the Monthly Trend sheet to "Monthly_PO_Qty".
# 2) Add a new "PO Forecast" worksheet after "Monthly Trend" containing the
#    forecast data (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1) Rename headers on existing sheets -----------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$newSheet.Cells.Item(2,1).Value = 44983.99999999999
$newSheet.Cells.Item(2,2).Value = 8
$newSheet.Cells.Item(2,3).Value = -26.02482169115098
$newSheet.Cells.Item(2,4).Value = 41.0868949193929
$newSheet.Cells.Item(3,1).Value = 45109.99999999999
$newSheet.Cells.Item(3,2).Value = 13
$newSheet.Cells.Item(3,3).Value = -20.93214706483073
$newSheet.Cells.Item(3,4).Value = 47.87292959729613
$newSheet.Cells.Item(4,1).Value = 45116.99999999999
$newSheet.Cells.Item(4,2).Value = 14
$newSheet.Cells.Item(4,3).Value = -20.51531272738129
$newSheet.Cells.Item(4,4).Value = 49.9604970250295
$newSheet.Cells.Item(5,1).Value = 45144.99999999999
$newSheet.Cells.Item(5,2).Value = 15
$newSheet.Cells.Item(5,3).Value = -18.51248436063039
$newSheet.Cells.Item(5,4).Value = 47.8979349625201
$newSheet.Cells.Item(6,1).Value = 45151.99999999999
$newSheet.Cells.Item(6,2).Value = 15
$newSheet.Cells.Item(6,3).Value = -16.38279692702327
$newSheet.Cells.Item(6,4).Value = 49.57061212692496
$newSheet.Cells.Item(7,1).Value = 45179.99999999999
$newSheet.Cells.Item(7,2).Value = 17
$newSheet.Cells.Item(7,3).Value = -16.83132614007526
$newSheet.Cells.Item(7,4).Value = 51.1682921175582
$newSheet.Cells.Item(8,1).Value = 45186.99999999999
$newSheet.Cells.Item(8,2).Value = 17
$newSheet.Cells.Item(8,3).Value = -17.15432604780744
$newSheet.Cells.Item(8,4).Value = 50.84305503946234
$newSheet.Cells.Item(9,1).Value = 45193.99999999999
$newSheet.Cells.Item(9,2).Value = 17
$newSheet.Cells.Item(9,3).Value = -16.42742560017087
$newSheet.Cells.Item(9,4).Value = 51.04465499450965
$newSheet.Cells.Item(10,1).Value = 45200.99999999999
$newSheet.Cells.Item(10,2).Value = 17
$newSheet.Cells.Item(10,3).Value = -17.12805791759456
$newSheet.Cells.Item(10,4).Value = 49.80170775823393
$newSheet.Cells.Item(11,1).Value = 45207.99999999999
$newSheet.Cells.Item(11,2).Value = 18
$newSheet.Cells.Item(11,3).Value = -14.79805399372206
$newSheet.Cells.Item(11,4).Value = 52.01032520852884
$newSheet.Cells.Item(12,1).Value = 45228.99999999999
$newSheet.Cells.Item(12,2).Value = 19
$newSheet.Cells.Item(12,3).Value = -14.0643306201989
$newSheet.Cells.Item(12,4).Value = 53.96051329376392
$newSheet.Cells.Item(13,1).Value = 45242.99999999999
$newSheet.Cells.Item(13,2).Value = 19
$newSheet.Cells.Item(13,3).Value = -16.12691789400182
$newSheet.Cells.Item(13,4).Value = 52.40342395371736
$newSheet.Cells.Item(14,1).Value = 45298.99999999999
$newSheet.Cells.Item(14,2).Value = 22
$newSheet.Cells.Item(14,3).Value = -10.49150438126315
$newSheet.Cells.Item(14,4).Value = 51.52804431266066
$newSheet.Cells.Item(15,1).Value = 45305.99999999999
$newSheet.Cells.Item(15,2).Value = 22
$newSheet.Cells.Item(15,3).Value = -11.25874855200596
$newSheet.Cells.Item(15,4).Value = 56.10642039304145
$newSheet.Cells.Item(16,1).Value = 45312.99999999999
$newSheet.Cells.Item(16,2).Value = 22
$newSheet.Cells.Item(16,3).Value = -11.12863592174762
$newSheet.Cells.Item(16,4).Value = 54.93134089140352
$newSheet.Cells.Item(17,1).Value = 45326.99999999999
$newSheet.Cells.Item(17,2).Value = 23
$newSheet.Cells.Item(17,3).Value = -9.136618404284214
$newSheet.Cells.Item(17,4).Value = 57.29850748063774
$newSheet.Cells.Item(18,1).Value = 45333.99999999999
$newSheet.Cells.Item(18,2).Value = 23
$newSheet.Cells.Item(18,3).Value = -10.72859645254793
$newSheet.Cells.Item(18,4).Value = 57.04107582276254
$newSheet.Cells.Item(19,1).Value = 45354.99999999999
$newSheet.Cells.Item(19,2).Value = 24
$newSheet.Cells.Item(19,3).Value = -8.005703064301104
$newSheet.Cells.Item(19,4).Value = 58.45675421509043
$newSheet.Cells.Item(20,1).Value = 45361.99999999999
$newSheet.Cells.Item(20,2).Value = 25
$newSheet.Cells.Item(20,3).Value = -9.460672637985549
$newSheet.Cells.Item(20,4).Value = 61.26815452883562
$newSheet.Cells.Item(21,1).Value = 45368.99999999999
$newSheet.Cells.Item(21,2).Value = 25
$newSheet.Cells.Item(21,3).Value = -7.741476526673332
$newSheet.Cells.Item(21,4).Value = 55.52621905337507
$newSheet.Cells.Item(22,1).Value = 45382.99999999999
$newSheet.Cells.Item(22,2).Value = 26
$newSheet.Cells.Item(22,3).Value = -7.237106254028508
$newSheet.Cells.Item(22,4).Value = 59.8151355112071
$newSheet.Cells.Item(23,1).Value = 45389.99999999999
$newSheet.Cells.Item(23,2).Value = 26
$newSheet.Cells.Item(23,3).Value = -5.613250475764293
$newSheet.Cells.Item(23,4).Value = 59.74963256659445
$newSheet.Cells.Item(24,1).Value = 45396.99999999999
$newSheet.Cells.Item(24,2).Value = 26
$newSheet.Cells.Item(24,3).Value = -8.138475023128608
$newSheet.Cells.Item(24,4).Value = 59.44109398023988
$newSheet.Cells.Item(25,1).Value = 45403.99999999999
$newSheet.Cells.Item(25,2).Value = 26
$newSheet.Cells.Item(25,3).Value = -6.640659982352402
$newSheet.Cells.Item(25,4).Value = 59.8794319747832
$newSheet.Cells.Item(26,1).Value = 45424.99999999999
$newSheet.Cells.Item(26,2).Value = 27
$newSheet.Cells.Item(26,3).Value = -4.766954208902018
$newSheet.Cells.Item(26,4).Value = 59.11910732076301
$newSheet.Cells.Item(27,1).Value = 45431.99999999999
$newSheet.Cells.Item(27,2).Value = 28
$newSheet.Cells.Item(27,3).Value = -4.392892630013441
$newSheet.Cells.Item(27,4).Value = 61.06887783149156
$newSheet.Cells.Item(28,1).Value = 45438.99999999999
$newSheet.Cells.Item(28,2).Value = 28
$newSheet.Cells.Item(28,3).Value = -5.474776326952492
$newSheet.Cells.Item(28,4).Value = 60.20234059167834
$newSheet.Cells.Item(29,1).Value = 45445.99999999999
$newSheet.Cells.Item(29,2).Value = 28
$newSheet.Cells.Item(29,3).Value = -2.237846806449673
$newSheet.Cells.Item(29,4).Value = 60.65159109648719
$newSheet.Cells.Item(30,1).Value = 45452.99999999999
$newSheet.Cells.Item(30,2).Value = 29
$newSheet.Cells.Item(30,3).Value = -3.795805172950262
$newSheet.Cells.Item(30,4).Value = 61.44285222919476
$newSheet.Cells.Item(31,1).Value = 45459.99999999999
$newSheet.Cells.Item(31,2).Value = 29
$newSheet.Cells.Item(31,3).Value = -3.885291764280487
$newSheet.Cells.Item(31,4).Value = 65.89740717228884
$newSheet.Cells.Item(32,1).Value = 45466.99999999999
$newSheet.Cells.Item(32,2).Value = 29
$newSheet.Cells.Item(32,3).Value = -5.048575059759017
$newSheet.Cells.Item(32,4).Value = 63.52612839503093
$newSheet.Cells.Item(33,1).Value = 45473.99999999999
$newSheet.Cells.Item(33,2).Value = 30
$newSheet.Cells.Item(33,3).Value = -4.589851857238429
$newSheet.Cells.Item(33,4).Value = 61.53142057078736
$newSheet.Cells.Item(34,1).Value = 45480.99999999999
$newSheet.Cells.Item(34,2).Value = 30
$newSheet.Cells.Item(34,3).Value = -2.382454514180466
$newSheet.Cells.Item(34,4).Value = 63.00231845047437
$newSheet.Cells.Item(35,1).Value = 45487.99999999999
$newSheet.Cells.Item(35,2).Value = 30
$newSheet.Cells.Item(35,3).Value = -5.881341161837479
$newSheet.Cells.Item(35,4).Value = 62.88813655103773
$newSheet.Cells.Item(36,1).Value = 45494.99999999999
$newSheet.Cells.Item(36,2).Value = 30
$newSheet.Cells.Item(36,3).Value = -2.839557245458218
$newSheet.Cells.Item(36,4).Value = 63.67386544625935
$newSheet.Cells.Item(37,1).Value = 45501.99999999999
$newSheet.Cells.Item(37,2).Value = 31
$newSheet.Cells.Item(37,3).Value = -1.537144150605482
$newSheet.Cells.Item(37,4).Value = 62.56579908512821
$newSheet.Cells.Item(38,1).Value = 45508.99999999999
$newSheet.Cells.Item(38,2).Value = 31
$newSheet.Cells.Item(38,3).Value = -4.809992599392067
$newSheet.Cells.Item(38,4).Value = 64.55025589167779

# Apply the date/time number format to the "ds" column (matches style used
# for the Order Week / Order Month date columns on the other sheets).
$newSheet.Range("A2:A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
